$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Insert 3 new columns before column B.
# This shifts the existing B..L ("venue_id" .. "geometry") to E..O
# and updates the sheet dimension from A1:L33 to A1:O33.
$ws.Range("B1:D1").EntireColumn.Insert()

# The insert operation copies column A's per-row formatting into the
# freshly inserted B:D columns; strip that back off the data rows so the
# new data cells end up unstyled, matching the source columns they mirror.
$ws.Range("B2:D33").ClearFormats()

# Re-apply the bold / bordered / centered header styling (as used by the
# other header cells) to the three new header cells by copying the format
# from the neighboring "venue_id" header cell (now at E1).
$ws.Range("E1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)

# Label the three new header columns.
$ws.Range("B1").Value = "Unnamed: 0.2"
$ws.Range("C1").Value = "Unnamed: 0.1"
$ws.Range("D1").Value = "Unnamed: 0"

# These three new columns duplicate the original row-index column (A) for
# every data row, just like pandas' extra "Unnamed: 0*" index columns.
$ws.Range("A2:A33").Copy()
$ws.Range("B2:B33").PasteSpecial(-4163)
$ws.Range("A2:A33").Copy()
$ws.Range("C2:C33").PasteSpecial(-4163)
$ws.Range("A2:A33").Copy()
$ws.Range("D2:D33").PasteSpecial(-4163)
